$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("choices")

# --- Insert a new "grappe" (cluster) choice row: LOGOGNEGUE 3, belonging to
#     csps TORANDOUGOU, placed right after LOGOGNEGUE 2 (old row 25, now new row 25).
#     This pushes the rest of the grappe list (old rows 25-34) down by one.
$ws.Rows.Item(25).Insert()
$ws.Cells.Item(25, 1).Value = "grappe"
$ws.Cells.Item(25, 2).Value = "LOGOGNEGUE 3"
$ws.Cells.Item(25, 3).Value = "LOGOGNEGUE 3"
$ws.Cells.Item(25, 5).Value = "TORANDOUGOU"

# --- Insert the matching "nb_grappe" row (the numbered list used to size the
#     cluster count) right after LOGOGNEGUE 2's entry (old row 46, now row 47
#     after the previous insertion shifted everything down by one).
$ws.Rows.Item(47).Insert()
$ws.Cells.Item(47, 1).Value = "nb_grappe"
$ws.Cells.Item(47, 2).Value = 12
$ws.Cells.Item(47, 3).Value = 12
$ws.Cells.Item(47, 6).Value = "LOGOGNEGUE 3"

# --- Renumber the remaining nb_grappe rows (old counts 12-20 become 13-21).
for ($r = 48; $r -le 56; $r++) {
    $newCount = $r - 35
    $ws.Cells.Item($r, 2).Value = $newCount
    $ws.Cells.Item($r, 3).Value = $newCount
}

# Update the choices sheet's selection to the new insertion point, without
# leaving it as the active tab.
$ws.Range("A46:A47").Select()

# The workbook ends up with "settings" as the active/selected tab.
$ws3 = $wb.Worksheets.Item("settings")
$ws3.Activate()
